$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new donor row (row 4): Kaviya / chennai / 9012346789.
# The id-proof number is stored as text (like the existing id-proof
# entries in column C), so it needs the leading apostrophe to keep
# Excel from coercing the all-digit string into a Number.
$ws.Range("A4").Value = "Kaviya"
$ws.Range("B4").Value = "chennai"
$ws.Range("C4").Value = "'9012346789"
